# Apply updated cryptos list (price + 1h volume change) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.784.90"
$ws.Range("E2").Value = "  -0.33%  "

$ws.Range("D3").Value = "2.042.75"
$ws.Range("E3").Value = "  +0.28%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("E6").Value = "  -0.90%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.16"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -2.64%  "

$ws.Range("E10").Value = "  +3.28%  "

$ws.Range("E11").Value = "  -0.18%  "

$ws.Range("E12").Value = "  +0.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.765"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.46%  "

$ws.Range("D17").Value = "2.033.78"
$ws.Range("E17").Value = "  -0.63%  "

$ws.Range("D18").Value = "37.759.40"
$ws.Range("E18").Value = "  -0.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.59%  "

$ws.Range("E20").Value = "  -2.65%  "

$ws.Range("E21").Value = "  +0.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.64%  "

$ws.Range("E23").Value = "  +0.53%  "

$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.98%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.64%  "

$ws.Range("E28").Value = "  -0.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.05%  "

$ws.Range("E30").Value = "  -0.26%  "

$ws.Range("E31").Value = "  -0.97%  "

$ws.Range("E32").Value = "  +8.60%  "

$ws.Range("E33").Value = "  -1.33%  "

$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("E35").Value = "  -0.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.04%  "

$ws.Range("E38").Value = "  +6.22%  "

$ws.Range("E39").Value = "  -0.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.45%  "

$ws.Range("D41").Value = "1.530.04"
$ws.Range("E41").Value = "  -0.97%  "

$ws.Range("E42").Value = "  +0.73%  "

$ws.Range("E43").Value = "  -0.99%  "

$ws.Range("E45").Value = "  -1.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.04%  "

$ws.Range("E47").Value = "  -0.15%  "

$ws.Range("E48").Value = "  +0.05%  "

$ws.Range("E49").Value = "  -0.64%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.24%  "

$ws.Range("D51").Value = "2.235.25"
$ws.Range("E51").Value = "  +0.35%  "
